$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new QR-scanner log entry as row 42 (mirrors the existing rows,
# which store every column - including the numeric-looking Student ID - as
# plain text).
$ws.Range("A42").Formula = '="201987"'
$ws.Range("B42").Formula = '="Neurology"'
$ws.Range("C42").Formula = '="28/12/2025"'
$ws.Range("D42").Formula = '="10:53:14"'
$ws.Range("E42").Formula = '="Manual"'
$ws.Range("F42").Formula = '="emp17.farah.a.youssef@gmail.com"'

# Convert the helper formulas to plain text values (so cells come out as
# real text, same as the rest of the sheet, and not formulas/numbers).
$ws.Range("A42:F42").Copy()
$ws.Range("A42:F42").PasteSpecial(-4163)
